$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt16"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4491763333333333
$ws.Range("H2").Value = 1.347529
$ws.Range("I2").Value = 0.236099761434867
$ws.Range("J2").Value = 0.236099761434867
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1811433333333334
$ws.Range("N2").Value = 0.5434300000000001
$ws.Range("O2").Value = 0.0111261749556462
$ws.Range("P2").Value = 0.01112617495564619
$ws.Range("Q2").Value = 0.08136529827444446
$ws.Range("R2").Value = 0.7322876844700001
$ws.Range("S2").Value = 0.002626887252710658
$ws.Range("T2").Value = 0.002626887252710658

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt16"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4491763333333333
$ws.Range("H3").Value = 1.347529
$ws.Range("I3").Value = 0.236099761434867
$ws.Range("J3").Value = 0.236099761434867
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.42533133333333
$ws.Range("N3").Value = 40.275994
$ws.Range("O3").Value = 0.8246098959508241
$ws.Range("P3").Value = 0.8246098959508241
$ws.Range("Q3").Value = 6.030341102091778
$ws.Range("R3").Value = 54.273069918826
$ws.Range("S3").Value = 0.19469019971082
$ws.Range("T3").Value = 0.19469019971082

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt16"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4491763333333333
$ws.Range("H4").Value = 1.347529
$ws.Range("I4").Value = 0.236099761434867
$ws.Range("J4").Value = 0.236099761434867
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.659118666666667
$ws.Range("N4").Value = 7.977356
$ws.Range("O4").Value = 0.1633282272592126
$ws.Range("P4").Value = 0.1633282272592126
$ws.Range("Q4").Value = 1.194413172591556
$ws.Range("R4").Value = 10.749718553324
$ws.Range("S4").Value = 0.03856175549147983
$ws.Range("T4").Value = 0.03856175549147983

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt16"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4491763333333333
$ws.Range("H5").Value = 1.347529
$ws.Range("I5").Value = 0.236099761434867
$ws.Range("J5").Value = 0.236099761434867
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.015234
$ws.Range("N5").Value = 0.045702
$ws.Range("O5").Value = 0.0009357018343171013
$ws.Range("P5").Value = 0.0009357018343171013
$ws.Range("Q5").Value = 0.006842752262
$ws.Range("R5").Value = 0.061584770358
$ws.Range("S5").Value = 0.000220918979856435
$ws.Range("T5").Value = 0.000220918979856435

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt16"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.190874
$ws.Range("H6").Value = 3.572622
$ws.Range("I6").Value = 0.6259569937989885
$ws.Range("J6").Value = 0.6259569937989886
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1811433333333334
$ws.Range("N6").Value = 0.5434300000000001
$ws.Range("O6").Value = 0.0111261749556462
$ws.Range("P6").Value = 0.01112617495564619
$ws.Range("Q6").Value = 0.21571888594
$ws.Range("R6").Value = 1.94146997346
$ws.Range("S6").Value = 0.006964507027717888
$ws.Range("T6").Value = 0.006964507027717888

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt16"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.190874
$ws.Range("H7").Value = 3.572622
$ws.Range("I7").Value = 0.6259569937989885
$ws.Range("J7").Value = 0.6259569937989886
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.42533133333333
$ws.Range("N7").Value = 40.275994
$ws.Range("O7").Value = 0.8246098959508241
$ws.Range("P7").Value = 0.8246098959508241
$ws.Range("Q7").Value = 15.987878026252
$ws.Range("R7").Value = 143.890902236268
$ws.Range("S7").Value = 0.5161703315262746
$ws.Range("T7").Value = 0.5161703315262747

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt16"
$ws.Range("C8").Value = "Fzd2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.190874
$ws.Range("H8").Value = 3.572622
$ws.Range("I8").Value = 0.6259569937989885
$ws.Range("J8").Value = 0.6259569937989886
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.659118666666667
$ws.Range("N8").Value = 7.977356
$ws.Range("O8").Value = 0.1633282272592126
$ws.Range("P8").Value = 0.1633282272592126
$ws.Range("Q8").Value = 3.166675283048
$ws.Range("R8").Value = 28.500077547432
$ws.Range("S8").Value = 0.1022364461376947
$ws.Range("T8").Value = 0.1022364461376947

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt16"
$ws.Range("C9").Value = "Fzd2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.190874
$ws.Range("H9").Value = 3.572622
$ws.Range("I9").Value = 0.6259569937989885
$ws.Range("J9").Value = 0.6259569937989886
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.015234
$ws.Range("N9").Value = 0.045702
$ws.Range("O9").Value = 0.0009357018343171013
$ws.Range("P9").Value = 0.0009357018343171013
$ws.Range("Q9").Value = 0.018141774516
$ws.Range("R9").Value = 0.163275970644
$ws.Range("S9").Value = 0.000585709107301332
$ws.Range("T9").Value = 0.0005857091073013321

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt16"
$ws.Range("C10").Value = "Fzd2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.06510733333333334
$ws.Range("H10").Value = 0.195322
$ws.Range("I10").Value = 0.03422225243611164
$ws.Range("J10").Value = 0.03422225243611164
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1811433333333334
$ws.Range("N10").Value = 0.5434300000000001
$ws.Range("O10").Value = 0.0111261749556462
$ws.Range("P10").Value = 0.01112617495564619
$ws.Range("Q10").Value = 0.01179375938444445
$ws.Range("R10").Value = 0.10614383446
$ws.Range("S10").Value = 0.0003807627679804674
$ws.Range("T10").Value = 0.0003807627679804673

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Wnt16"
$ws.Range("C11").Value = "Fzd2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.06510733333333334
$ws.Range("H11").Value = 0.195322
$ws.Range("I11").Value = 0.03422225243611164
$ws.Range("J11").Value = 0.03422225243611164
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.42533133333333
$ws.Range("N11").Value = 40.275994
$ws.Range("O11").Value = 0.8246098959508241
$ws.Range("P11").Value = 0.8246098959508241
$ws.Range("Q11").Value = 0.8740875222297778
$ws.Range("R11").Value = 7.866787700068
$ws.Range("S11").Value = 0.02822000802054486
$ws.Range("T11").Value = 0.02822000802054486

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Wnt16"
$ws.Range("C12").Value = "Fzd2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.06510733333333334
$ws.Range("H12").Value = 0.195322
$ws.Range("I12").Value = 0.03422225243611164
$ws.Range("J12").Value = 0.03422225243611164
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.659118666666667
$ws.Range("N12").Value = 7.977356
$ws.Range("O12").Value = 0.1633282272592126
$ws.Range("P12").Value = 0.1633282272592126
$ws.Range("Q12").Value = 0.1731281254035556
$ws.Range("R12").Value = 1.558153128632
$ws.Range("S12").Value = 0.005589459823207385
$ws.Range("T12").Value = 0.005589459823207384

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Wnt16"
$ws.Range("C13").Value = "Fzd2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.06510733333333334
$ws.Range("H13").Value = 0.195322
$ws.Range("I13").Value = 0.03422225243611164
$ws.Range("J13").Value = 0.03422225243611164
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.015234
$ws.Range("N13").Value = 0.045702
$ws.Range("O13").Value = 0.0009357018343171013
$ws.Range("P13").Value = 0.0009357018343171013
$ws.Range("Q13").Value = 0.000991845116
$ws.Range("R13").Value = 0.008926606044
$ws.Range("S13").Value = 0.00003202182437893255
$ws.Range("T13").Value = 0.00003202182437893255

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Wnt16"
$ws.Range("C14").Value = "Fzd2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1973276666666667
$ws.Range("H14").Value = 0.5919829999999999
$ws.Range("I14").Value = 0.1037209923300328
$ws.Range("J14").Value = 0.1037209923300328
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.1811433333333334
$ws.Range("N14").Value = 0.5434300000000001
$ws.Range("O14").Value = 0.0111261749556462
$ws.Range("P14").Value = 0.01112617495564619
$ws.Range("Q14").Value = 0.03574459129888889
$ws.Range("R14").Value = 0.32170132169
$ws.Range("S14").Value = 0.001154017907237183
$ws.Range("T14").Value = 0.001154017907237183

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Wnt16"
$ws.Range("C15").Value = "Fzd2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1973276666666667
$ws.Range("H15").Value = 0.5919829999999999
$ws.Range("I15").Value = 0.1037209923300328
$ws.Range("J15").Value = 0.1037209923300328
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 13.42533133333333
$ws.Range("N15").Value = 40.275994
$ws.Range("O15").Value = 0.8246098959508241
$ws.Range("P15").Value = 0.8246098959508241
$ws.Range("Q15").Value = 2.649189306233555
$ws.Range("R15").Value = 23.84270375610199
$ws.Range("S15").Value = 0.08552935669318461
$ws.Range("T15").Value = 0.08552935669318461

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Wnt16"
$ws.Range("C16").Value = "Fzd2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1973276666666667
$ws.Range("H16").Value = 0.5919829999999999
$ws.Range("I16").Value = 0.1037209923300328
$ws.Range("J16").Value = 0.1037209923300328
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.659118666666667
$ws.Range("N16").Value = 7.977356
$ws.Range("O16").Value = 0.1633282272592126
$ws.Range("P16").Value = 0.1633282272592126
$ws.Range("Q16").Value = 0.5247176818831111
$ws.Range("R16").Value = 4.722459136947999
$ws.Range("S16").Value = 0.01694056580683065
$ws.Range("T16").Value = 0.01694056580683065

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Wnt16"
$ws.Range("C17").Value = "Fzd2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1973276666666667
$ws.Range("H17").Value = 0.5919829999999999
$ws.Range("I17").Value = 0.1037209923300328
$ws.Range("J17").Value = 0.1037209923300328
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.015234
$ws.Range("N17").Value = 0.045702
$ws.Range("O17").Value = 0.0009357018343171013
$ws.Range("P17").Value = 0.0009357018343171013
$ws.Range("Q17").Value = 0.003006089673999999
$ws.Range("R17").Value = 0.027054807066
$ws.Range("S17").Value = 0.00009705192278040174
$ws.Range("T17").Value = 0.00009705192278040174
